$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7333130836486816
$ws.Range("B1").Value = 2.931930780410767
$ws.Range("C1").Value = 3.093540906906128
$ws.Range("D1").Value = 2.450787305831909
$ws.Range("E1").Value = 1.499755024909973
